$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update B2 and add E2 formula ---
$ws.Range("B2").Value = 1792
$ws.Range("E2").Formula = "=D2/2"

# --- Row 3: update B3 and add E3 formula ---
$ws.Range("B3").Value = 2560
$ws.Range("E3").Formula = "=D3/2"

# --- Row 5: add E5 value (matches B5/C5/D5) ---
$ws.Range("E5").Value = 32

# --- Row 8: add E8 formula ---
$ws.Range("E8").Formula = "=E2/E5"

# --- Row 9: add E9 formula ---
$ws.Range("E9").Formula = "=E3/E5"

# --- New rows 14-15 ---
$ws.Range("B14").Value = 720
$ws.Range("B15").Value = 1280

# --- New rows 18-19: spiral/block search calculations ---
$ws.Range("B18").Value = 584
$ws.Range("C18").Formula = "=B18*4"
$ws.Range("E18").Formula = "=B3-C18"
$ws.Range("F18").Formula = "=E18/2"

$ws.Range("B19").Value = 388
$ws.Range("C19").Formula = "=B19*4"
$ws.Range("E19").Formula = "=B2-C19"
$ws.Range("F19").Formula = "=E19/2"

# --- Update selection to match target state ---
$ws.Range("E18").Select()
